$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.620.86"
$ws.Range("E2").Value = "  -0.06%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.595.94"
$ws.Range("E3").Value = "  +0.44%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.03%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'211.55"
$ws.Range("E5").Value = "  +0.30%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +1.32%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.01%  "

# Row 8 - Dogecoin
$ws.Range("D8").Value = "'0.0616"

# Row 9 - Cardano
$ws.Range("E9").Value = "  -0.46%  "

# Row 10 - Solana
$ws.Range("D10").Value = "'19.42"
$ws.Range("E10").Value = "  -0.80%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.60%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.819.68"
$ws.Range("E12").Value = "  +0.47%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.578.48"
$ws.Range("E13").Value = "  -0.64%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  +0.27%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  +0.15%  "

# Row 16 - Litecoin
$ws.Range("D16").Value = "'64.55"
$ws.Range("E16").Value = "  -0.20%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "26.600.44"
$ws.Range("E17").Value = "  -0.04%  "

# Row 18 - ShibaInu
$ws.Range("E18").Value = "  +0.40%  "

# Row 19 - BitcoinCash
$ws.Range("D19").Value = "'208.52"
$ws.Range("E19").Value = "  +0.28%  "

# Row 20 - Dai
$ws.Range("E20").Value = "  -0.04%  "

# Row 21 - Chainlink
$ws.Range("E21").Value = "  +3.51%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  +0.40%  "

# Row 23 - Toncoin
$ws.Range("E23").Value = "  -1.51%  "

# Row 24 - Avalanche
$ws.Range("D24").Value = "'8.87"
$ws.Range("E24").Value = "  +0.10%  "

# Row 25 - Monero
$ws.Range("E25").Value = "  -1.13%  "

# Row 26 - BinanceUSD
$ws.Range("E26").Value = "  +0.06%  "

# Row 27 - Cosmos
$ws.Range("E27").Value = "  -1.95%  "

# Row 28 - Stellar
$ws.Range("D28").Value = "'0.114"
$ws.Range("E28").Value = "  +0.92%  "

# Row 29 - EthereumClassic
$ws.Range("D29").Value = "'15.23"
$ws.Range("E29").Value = "  -0.16%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  -0.07%  "

# Row 31 - PancakeSwap
$ws.Range("D31").Value = "'1.16"
$ws.Range("E31").Value = "  +0.55%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  +0.25%  "

# Rows 33 and 34 - swap InternetComputer(DFINITY) and ImmutableX, with updated values
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "'0.652"
$ws.Range("E33").Value = "  -1.33%  "

$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "'2.93"
$ws.Range("E34").Value = "  +0.92%  "

# Row 35 - Maker
$ws.Range("D35").Value = "1.281.26"
$ws.Range("E35").Value = "  -1.85%  "

# Row 36 - HuobiToken
$ws.Range("E36").Value = "  +0.90%  "

# Row 38 - VeChain
$ws.Range("E38").Value = "  -0.32%  "

# Row 39 - ARBITRUM
$ws.Range("E39").Value = "  +1.78%  "

# Row 41 - FraxShare
$ws.Range("E41").Value = "  +2.18%  "

# Row 42 - MXToken
$ws.Range("E42").Value = "  +1.26%  "

# Row 43 - TrustWalletToken
$ws.Range("E43").Value = "  -0.80%  "

# Row 44 - Aave
$ws.Range("D44").Value = "'64.08"
$ws.Range("E44").Value = "  +2.41%  "

# Row 45 - RocketPoolETH
$ws.Range("D45").Value = "1.731.65"
$ws.Range("E45").Value = "  +0.44%  "

# Row 46 - WEMIXToken
$ws.Range("D46").Value = "'0.910"
$ws.Range("E46").Value = "  +8.58%  "

# Row 47 - Quant
$ws.Range("D47").Value = "'89.58"
$ws.Range("E47").Value = "  +0.11%  "

# Row 48 - RenderToken
$ws.Range("E48").Value = "  -0.80%  "

# Row 49 - BabyDogeCoin
$ws.Range("E49").Value = "  -2.10%  "

# Row 50 - Algorand
$ws.Range("E50").Value = "  +4.41%  "

# Row 51 - Cronos
$ws.Range("E51").Value = "  +0.40%  "
